$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '43.741.15'
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '  -0.72%  '
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.227.24'
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '  -1.17%  '
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '  +0.19%  '
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '272.47'
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '  +5.51%  '
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '86.76'
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '  +9.26%  '
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.621'
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '  -0.96%  '
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '  +0.09%  '
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.603'
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '  -0.72%  '
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '44.79'
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '  +3.71%  '
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '  -1.03%  '
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '7.68'
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '  +7.99%  '
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '  +1.07%  '
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '2.570.57'
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '  -0.77%  '
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '14.84'
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '  +0.70%  '
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '2.227.90'
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '  -0.44%  '
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.789'
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '  -0.78%  '
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '43.713.95'
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '  -0.57%  '
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '  -1.62%  '
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '70.18'
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '5.94'
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '  -2.27%  '
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '  -0.67%  '
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '232.41'
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '  -1.12%  '
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '  -9.73%  '
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '  +0.08%  '
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.53'
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '  +13.14%  '
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '3.54'
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '  +4.88%  '
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '10.74'
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '  -1.03%  '
$cell.Style = "Normal"

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = '  +5.36%  '
$cell.Style = "Normal"

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '38.88'
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = '  -8.26%  '
$cell.Style = "Normal"

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '173.42'
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = '  -0.16%  '
$cell.Style = "Normal"

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.0899'
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = '  +2.11%  '
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '20.69'
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = '  -0.07%  '
$cell.Style = "Normal"

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '5.33'
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = '  +0.14%  '
$cell.Style = "Normal"

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = '  -0.14%  '
$cell.Style = "Normal"

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = '  -2.58%  '
$cell.Style = "Normal"

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.0350'
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = '  -3.14%  '
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '4.29'
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '  -4.88%  '
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '  +16.27%  '
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.19'
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '  +2.01%  '
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '12.36'
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '  -8.87%  '
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '63.39'
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '  +1.21%  '
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '  -0.19%  '
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '5.38'
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '  -0.79%  '
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '8.48'
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '  -1.42%  '
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.0988'
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '  -0.35%  '
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '99.70'
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '  -5.46%  '
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '  +3.52%  '
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '  +0.14%  '
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '  -2.34%  '
$cell.Style = "Normal"

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '  -10.12%  '
$cell.Style = "Normal"
